$wb = $excel.ActiveWorkbook

# 1. Rename Sheet2 -> Paxdetails
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Paxdetails"

# 2a. Main pax header: "Pax 1" -> "Mainpax", bold
$ws.Range("A1").Value = "Mainpax"
$ws.Range("A1").Font.Bold = $true

# 2b. Child pax header at row 9
$ws.Range("A9").Value = "Childpax"
$ws.Range("A9").Font.Bold = $true

# 2c. Infant pax header at row 18
$ws.Range("A18").Value = "Infantpax"
$ws.Range("A18").Font.Bold = $true

# 2d. Infant firstname change
$ws.Range("B21").Value = "Parvin"

# 2e. New 2nd-adult pax block, rows 27-31
$ws.Range("A27").Value = "2Adultpax"
$ws.Range("A27").Font.Bold = $true

$ws.Range("A28").Value = "Title"
$ws.Range("B28").Value = "Mr."

$ws.Range("A29").Value = "Firstname"
$ws.Range("B29").Value = "Pradeep"

$ws.Range("A30").Value = "Lastname"
$ws.Range("B30").Value = "GQ"

$ws.Range("A31").Value = "Residential"
$ws.Range("B31").Value = "2000-BARANGAROO,NSW"

# 3. Select whole row 28 and activate the Paxdetails sheet
$ws.Activate()
$ws.Rows.Item(28).Select()
